$wb = $excel.ActiveWorkbook

# --- "Cronograma #1" sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Cronograma #1")
$ws.Activate()

# Progress column (H) for the finished Sprint #1 tasks: 0% -> 100%
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
$ws.Range("H22").Value = 1
$ws.Range("H23").Value = 1

# Row 37 no longer needs the taller wrapped height
$ws.Rows.Item(37).RowHeight = 14.25

# Update the sheet's saved scroll position / selection
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("B37").Select()

# --- "IntegrantesRoles" sheet -----------------------------------------------
$ws2 = $wb.Worksheets.Item("IntegrantesRoles")
for ($i = 5; $i -le 11; $i++) {
    $ws2.Rows.Item($i).RowHeight = 12.75
}
